$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 287, shifting existing rows 287-344 down to 288-345
$ws.Rows.Item(287).Insert()

# Populate the newly inserted row 287 with the new data record
$ws.Range("A287").Value = 10
$ws.Range("B287").Value = "Vega Modelo de Temuco"
$ws.Range("C287").Value = "La Araucanía"
$ws.Range("D287").Value = 44816
$ws.Range("E287").Value = 9
$ws.Range("F287").Value = 100112001
$ws.Range("G287").Value = "Berenjena"
$ws.Range("H287").Value = "Sin especificar"
$ws.Range("I287").Value = "Primera"
$ws.Range("J287").Value = 85
$ws.Range("K287").Value = 15000
$ws.Range("L287").Value = 16000
$ws.Range("M287").Value = 15353
$ws.Range("N287").Value = '$/caja 40 unidades'
$ws.Range("O287").Value = "Región de Arica y Parinacota"
$ws.Range("P287").Value = 384
$ws.Range("Q287").Value = 40
$ws.Range("R287").Value = "Hortaliza"
